$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -10
$ws.Range("F7").Value = -8
$ws.Range("F11").Value = -9
$ws.Range("F12").Value = -7
$ws.Range("F13").Value = -3
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = -3
$ws.Range("F18").Value = -7
$ws.Range("F20").Value = -8
$ws.Range("F22").Value = -8
$ws.Range("F23").Value = -1
$ws.Range("F25").Value = 5
$ws.Range("F27").Value = 4
$ws.Range("F30").Value = -2

$wb.Save()
